$wb = $excel.ActiveWorkbook

# The file "e75c870a-5423-494b-851c-b458ce5c89e8.md" has been handed off again
# (it is now "Ready for handoff" instead of previously handed-back), so update
# its Status / Latest Handoff Datetime across the Overview, zh-cn, and de-de
# sheets.

# --- Overview sheet (row 3 = e75c870a-5423-494b-851c-b458ce5c89e8.md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-18 16:54:07"

# --- zh-cn sheet (row 3 = e75c870a-5423-494b-851c-b458ce5c89e8.md) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-18 16:53:58"

# --- de-de sheet (row 3 = e75c870a-5423-494b-851c-b458ce5c89e8.md) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-18 16:54:07"
